$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear old layout remnants in row 1 (B1:D1) - data moves down to rows 2-3
$ws.Range("B1:D1").Clear()

# Row 1: keep A1 = performLogin (unchanged)

# Row 2: pune / indore (unchanged content, just shifted down a row)
$ws.Range("B2").Value = "pune"
$ws.Range("C2").Value = "indore"

# Row 3: Indore / bhopal
$ws.Range("B3").Value = "Indore"
$ws.Range("C3").Value = "bhopal"

# Dates (text-formatted), entered last
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "07-Apr-2021"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "07-Mar-2021"

[void]$ws.Range("A1").Select()
